$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 110.68733978271484
$ws.Range("B3").Value = 107.13040924072266
$ws.Range("B4").Value = 104.14569091796875
$ws.Range("B5").Value = 101.14502716064453
$ws.Range("B6").Value = 101.46086120605469
$ws.Range("B7").Value = 101.52155303955078
$ws.Range("B8").Value = 99.152023315429688
$ws.Range("B9").Value = 99.53228759765625
$ws.Range("B10").Value = 124.6224365234375
$ws.Range("B11").Value = 149.14256286621094
$ws.Range("B12").Value = 116.40653991699219

$wb.Save()
